$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.45
$ws.Range("I3").Value = 3
$ws.Range("N3").Value = 2.63
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 1.57
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 2.1
$ws.Range("S3").Value = 1.67
$ws.Range("AC3").Value = 67
$ws.Range("AI3").Value = 29

# Row 5
$ws.Range("G5").Value = 1.65
$ws.Range("H5").Value = 4.1
$ws.Range("I5").Value = 4.75
$ws.Range("L5").Value = 1.22
$ws.Range("M5").Value = 4.33
$ws.Range("N5").Value = 1.73
$ws.Range("O5").Value = 2.1
$ws.Range("AA5").Value = 8
$ws.Range("AF5").Value = 26

# Row 7
$ws.Range("G7").Value = 1.28
$ws.Range("H7").Value = 4.6
$ws.Range("I7").Value = 8.75
$ws.Range("L7").Value = 1.17
$ws.Range("M7").Value = 4.4
$ws.Range("N7").Value = 1.53
$ws.Range("O7").Value = 2.18
$ws.Range("R7").Value = 1.94
$ws.Range("S7").Value = 1.77
$ws.Range("T7").Value = 6.5
$ws.Range("V7").Value = 7.2
$ws.Range("W7").Value = 7
$ws.Range("Y7").Value = 20
$ws.Range("Z7").Value = 13
$ws.Range("AA7").Value = 8.25
$ws.Range("AB7").Value = 17
$ws.Range("AC7").Value = 70
$ws.Range("AD7").Value = 450
$ws.Range("AE7").Value = 19.5
$ws.Range("AF7").Value = 50
$ws.Range("AG7").Value = 22
$ws.Range("AH7").Value = 200
$ws.Range("AI7").Value = 80
$ws.Range("AJ7").Value = 65

# Row 9
$ws.Range("H9").Value = 3.7
$ws.Range("I9").Value = 2.9
$ws.Range("J9").Value = 1.03
$ws.Range("K9").Value = 15
$ws.Range("L9").Value = 1.22
$ws.Range("M9").Value = 4
$ws.Range("N9").Value = 1.73
$ws.Range("O9").Value = 2.08
$ws.Range("P9").Value = 1.33
$ws.Range("Q9").Value = 3.25
$ws.Range("R9").Value = 1.67
$ws.Range("S9").Value = 2.1
$ws.Range("T9").Value = 9
$ws.Range("Y9").Value = 23
$ws.Range("Z9").Value = 13
$ws.Range("AA9").Value = 7.5
$ws.Range("AB9").Value = 13
$ws.Range("AC9").Value = 41
$ws.Range("AD9").Value = 151
$ws.Range("AE9").Value = 11

# Row 10
$ws.Range("G10").Value = 1.73
$ws.Range("H10").Value = 3.9
$ws.Range("I10").Value = 3.7
$ws.Range("R10").Value = 1.5
$ws.Range("S10").Value = 2.5
$ws.Range("U10").Value = 11
$ws.Range("X10").Value = 13
$ws.Range("AC10").Value = 34
$ws.Range("AD10").Value = 101
$ws.Range("AF10").Value = 23
$ws.Range("AG10").Value = 13
$ws.Range("AI10").Value = 26

# Row 11
$ws.Range("I11").Value = 3.9
$ws.Range("U11").Value = 9.5
$ws.Range("AE11").Value = 13
$ws.Range("AF11").Value = 21
$ws.Range("AG11").Value = 13
$ws.Range("AI11").Value = 29

# Row 12
$ws.Range("G12").Value = 2.15
$ws.Range("H12").Value = 3.25
$ws.Range("J12").Value = 1.05
$ws.Range("K12").Value = 11
$ws.Range("N12").Value = 1.93
$ws.Range("O12").Value = 1.93
$ws.Range("T12").Value = 8.5
$ws.Range("Z12").Value = 11
$ws.Range("AE12").Value = 11

# Row 13
$ws.Range("G13").Value = 2.2
$ws.Range("H13").Value = 3.05
$ws.Range("I13").Value = 3.15
$ws.Range("X13").Value = 19
$ws.Range("AA13").Value = 6
$ws.Range("AB13").Value = 15
$ws.Range("AE13").Value = 8.5
$ws.Range("AF13").Value = 16
$ws.Range("AG13").Value = 11.5
$ws.Range("AH13").Value = 45
$ws.Range("AI13").Value = 32

# Row 15
$ws.Range("G15").Value = 1.65
$ws.Range("H15").Value = 3.55
$ws.Range("I15").Value = 4.45
$ws.Range("N15").Value = 1.75
$ws.Range("O15").Value = 1.87
$ws.Range("R15").Value = 1.76
$ws.Range("S15").Value = 1.96
$ws.Range("T15").Value = 6.2
$ws.Range("U15").Value = 7
$ws.Range("V15").Value = 6.9
$ws.Range("W15").Value = 11
$ws.Range("X15").Value = 10.75
$ws.Range("Y15").Value = 19
$ws.Range("Z15").Value = 10.5
$ws.Range("AA15").Value = 6.1
$ws.Range("AB15").Value = 12
$ws.Range("AC15").Value = 50
$ws.Range("AE15").Value = 10.75
$ws.Range("AF15").Value = 21
$ws.Range("AG15").Value = 12
$ws.Range("AH15").Value = 55
$ws.Range("AI15").Value = 32
$ws.Range("AJ15").Value = 32

# Row 16
$ws.Range("G16").Value = 1.09
$ws.Range("H16").Value = 7.6
$ws.Range("I16").Value = 15.5
$ws.Range("O16").Value = 3.45
$ws.Range("R16").Value = 2.15
$ws.Range("S16").Value = 1.62
$ws.Range("U16").Value = 6.3
$ws.Range("X16").Value = 9.5
$ws.Range("Y16").Value = 28
$ws.Range("Z16").Value = 23
$ws.Range("AA16").Value = 16.5
$ws.Range("AB16").Value = 30
$ws.Range("AE16").Value = 45
$ws.Range("AF16").Value = 150
$ws.Range("AG16").Value = 45
$ws.Range("AH16").Value = 600
$ws.Range("AI16").Value = 200
$ws.Range("AJ16").Value = 120

# Row 17
$ws.Range("G17").Value = 2.25
$ws.Range("H17").Value = 3.2
$ws.Range("I17").Value = 3.3
$ws.Range("K17").Value = 8.5
$ws.Range("P17").Value = 1.44
$ws.Range("Q17").Value = 2.63
$ws.Range("R17").Value = 1.91
$ws.Range("S17").Value = 1.91
$ws.Range("X17").Value = 19
$ws.Range("Y17").Value = 29
$ws.Range("Z17").Value = 8.5
$ws.Range("AD17").Value = 301
$ws.Range("AE17").Value = 9

# Row 18
$ws.Range("G18").Value = 1.29
$ws.Range("I18").Value = 9
$ws.Range("J18").Value = 1.03
$ws.Range("K18").Value = 17
$ws.Range("N18").Value = 1.5
$ws.Range("O18").Value = 2.5
$ws.Range("P18").Value = 1.25
$ws.Range("AD18").Value = 251

# Row 20
$ws.Range("J20").Value = 1.07
$ws.Range("K20").Value = 9
$ws.Range("N20").Value = 2.05
$ws.Range("O20").Value = 1.8
$ws.Range("P20").Value = 1.44
$ws.Range("Q20").Value = 2.63
$ws.Range("Z20").Value = 9

# Row 22
$ws.Range("H22").Value = 4.15
$ws.Range("I22").Value = 6.6
$ws.Range("K22").Value = 7.7
$ws.Range("L22").Value = 1.26
$ws.Range("M22").Value = 3.45
$ws.Range("N22").Value = 1.78
$ws.Range("O22").Value = 1.93
$ws.Range("P22").Value = 1.38
$ws.Range("Q22").Value = 2.8
$ws.Range("R22").Value = 1.95
$ws.Range("T22").Value = 6.6
$ws.Range("U22").Value = 6.6
$ws.Range("X22").Value = 11.75
$ws.Range("Y22").Value = 28
$ws.Range("Z22").Value = 7.7
$ws.Range("AA22").Value = 8.25
$ws.Range("AD22").Value = 800
$ws.Range("AI22").Value = 75
